$d = $word.ActiveDocument

# The document's first paragraph currently reads (runs, then bookmark):
#   <w:r>Python adatelemzés</w:r><w:r> – SQL és Dataframe</w:r>
#   <w:bookmarkStart w:name="_GoBack"/><w:bookmarkEnd/>
# The target state moves the (collapsed/empty) "_GoBack" bookmark to the
# front of that same paragraph, before the two runs, leaving the runs
# themselves untouched:
#   <w:bookmarkStart w:name="_GoBack"/><w:bookmarkEnd/>
#   <w:r>Python adatelemzés</w:r><w:r> – SQL és Dataframe</w:r>

$target = $d.Paragraphs(1)
$destPos = $target.Range.Start

# Drop the existing bookmark from its current location (end of the
# paragraph's text, right before the paragraph mark).
$bk = $d.Bookmarks.Item("_GoBack")
$bk.Delete()

if ($destPos -eq 0) {
    # Creating a brand-new *collapsed* (empty) bookmark exactly at
    # absolute document position 0 is a degenerate case for this
    # object model (it silently expands to cover the whole paragraph
    # instead of staying collapsed), so work around it: temporarily
    # insert a throwaway empty paragraph right before the target
    # paragraph, anchor the new bookmark at the (now non-zero) start
    # of the target paragraph, then remove the throwaway paragraph
    # again - the bookmark (and the target paragraph's runs) shift
    # back down to position 0 intact.
    $target.Range.InsertParagraphBefore()

    $shiftedStart = $d.Paragraphs(2).Range.Start
    $anchor = $d.Range($shiftedStart, $shiftedStart)
    $d.Bookmarks.Add("_GoBack", $anchor)

    $d.Paragraphs(1).Range.Delete()
} else {
    $anchor = $d.Range($destPos, $destPos)
    $d.Bookmarks.Add("_GoBack", $anchor)
}

Write-Output "done"
